# Weekly data refresh: a new price record for "Cebollín" (Vega Modelo de
# Temuco) is inserted as row 274, pushing every subsequent record down by
# one row (the table keeps growing by one new weekly entry at the top of
# this date-cluster, while the previously-last row 335 slides down to the
# new row 336).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 274; Excel shifts rows 274:335 down to 275:336
# and keeps the existing column D date-number format on the new row.
$ws.Rows.Item(274).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(274, 1).Value = 10
$ws.Cells.Item(274, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(274, 3).Value = "La Araucanía"
$ws.Cells.Item(274, 4).Value = 44642
$ws.Cells.Item(274, 5).Value = 9
$ws.Cells.Item(274, 6).Value = 100112037
$ws.Cells.Item(274, 7).Value = "Cebollín"
$ws.Cells.Item(274, 8).Value = "Sin especificar"
$ws.Cells.Item(274, 9).Value = "Primera"
$ws.Cells.Item(274, 10).Value = 55
$ws.Cells.Item(274, 11).Value = 8000
$ws.Cells.Item(274, 12).Value = 8000
$ws.Cells.Item(274, 13).Value = 8000
$ws.Cells.Item(274, 14).Value = "`$/docena de paquetes"
$ws.Cells.Item(274, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(274, 16).Value = 667
$ws.Cells.Item(274, 17).Value = 12
$ws.Cells.Item(274, 18).Value = "Hortaliza"
